$wb = $excel.ActiveWorkbook

# --- Sheet: Aggr_generation ---
$wsGen = $wb.Worksheets.Item("Aggr_generation")
$wsGen.Range("B2").Value = 6066
$wsGen.Range("B3").Value = 1104
$wsGen.Range("B4").Value = 3768
$wsGen.Range("B5").Value = 1781
$wsGen.Range("B6").Value = 3175
$wsGen.Range("B7").Value = 1939
$wsGen.Range("B8").Value = 2468
$wsGen.Range("B9").Value = 4905
$wsGen.Range("B10").Value = 7280
$wsGen.Range("B11").Value = 1638

# --- Sheet: aggr_exchange ---
$wsExch = $wb.Worksheets.Item("aggr_exchange")
$wsExch.Range("B5").Value = 0
$wsExch.Range("B6").Value = -700
$wsExch.Range("B11").Value = -700
